# distintosPorcentajesLsPareto.xlsx — "excel actualizados a GKD only"
#
# The 5 non-GKD instance rows (MDG-a, MDG-b, SOM-a x2, SOM-b) that used to
# live in rows 34-38 are removed from the results table. The summary rows
# (3-12) use AVERAGE() formulas over B15:B38 etc., so clearing out those
# five rows' data shifts every averaged KPI; Excel recalculates
# automatically. The now-unused shared strings for those instance names
# drop out of the shared string table too, which is why every string index
# >= 105 in the original file shifts down by 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the data for the 5 removed instances (MDG-a_9, MDG-b_12,
# SOM-a_18, SOM-a_33, SOM-b_2) — rows 34 through 38, all 71 columns
# (A:BS). ClearContents keeps the cell's style (s="1") but drops the
# value/formula/shared-string reference, matching the target where column
# A's label cell disappears entirely (it had no explicit style) while
# B:BS keep their empty, styled cells.
$ws.Range("A34:BS38").ClearContents()

# The sheet was scrolled/selected elsewhere when last saved (topLeftCell
# AJ1, selection H24). Restore the default top-left scroll position and
# move the selection to A26.
[void]$ws.Range("A26").Select()
